# "add to dynamic grid" - fill in the Name/Email header and complete the
# fewest-coins dynamic-programming table (rows for N, D, Q denominations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: identify the student.
$ws.Range("A1").Value = "Name: Leo Liang"
$ws.Range("A2").Value = "Email: liang.jiahao@northeastern.edu"

# Row 15 (N = 5 cent coin): finish the remaining columns (40,45,50 cents).
$ws.Range("G15").Value = "4N"
$ws.Range("H15").Value = "5N"
$ws.Range("I15").Value = "6N"
$ws.Range("J15").Value = "7N"
$ws.Range("K15").Value = "8N"
$ws.Range("L15").Value = "9N"
$ws.Range("M15").Value = "10N"

# Row 16 (D = 10 cent coin): fill in the fewest-coin combination for each value.
$ws.Range("D16").Value = "1N"
$ws.Range("E16").Value = "1D"
$ws.Range("F16").Value = "1D 1N"
$ws.Range("G16").Value = "2D"
$ws.Range("H16").Value = "2D 1N"
$ws.Range("I16").Value = "3D"
$ws.Range("J16").Value = "3D 1N"
$ws.Range("K16").Value = "4D"
$ws.Range("L16").Value = "4D 1N"
$ws.Range("M16").Value = "5D"

# Row 17 (Q = 25 cent coin): fill in the fewest-coin combination for each value.
$ws.Range("D17").Value = "1N"
$ws.Range("E17").Value = "1D"
$ws.Range("F17").Value = "1D 1N"
$ws.Range("G17").Value = "2D"
$ws.Range("H17").Value = "1Q"
$ws.Range("I17").Value = "1Q 1N"
$ws.Range("J17").Value = "1Q 1D"
$ws.Range("K17").Value = "1Q 1D 1N"
$ws.Range("L17").Value = "1Q 2D"
$ws.Range("M17").Value = "2Q"

# Restore cursor/selection to match the finished worksheet.
$ws.Range("A12").Select()
